$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.341.46'
$ws.Cells.Item(2, 5).Value = '  +1.39%  '
$ws.Cells.Item(3, 4).Value = '1.828.43'
$ws.Cells.Item(3, 5).Value = '  +2.93%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  -0.18%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '317.81'
$ws.Cells.Item(5, 5).Value = '  +0.70%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.002'
$ws.Cells.Item(6, 5).Value = '  -0.22%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5344'
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4070'
$ws.Cells.Item(8, 5).Value = '  +9.33%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07610'
$ws.Cells.Item(9, 5).Value = '  +2.89%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '41.90'
$ws.Cells.Item(10, 5).Value = '  +0.93%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.104'
$ws.Cells.Item(11, 5).Value = '  +1.35%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '6.341'
$ws.Cells.Item(12, 5).Value = '  +4.50%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.003'
$ws.Cells.Item(13, 5).Value = '  -0.13%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.573'
$ws.Cells.Item(14, 5).Value = '  +5.27%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '20.81'
$ws.Cells.Item(15, 5).Value = '  +2.03%  '
$ws.Cells.Item(16, 4).Value = '1.828.08'
$ws.Cells.Item(16, 5).Value = '  +2.85%  '
$ws.Cells.Item(17, 2).Value = 'Litecoin'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '89.33'
$ws.Cells.Item(17, 5).Value = '  +1.31%  '
$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001074'
$ws.Cells.Item(18, 5).Value = '  +2.39%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06621'
$ws.Cells.Item(19, 5).Value = '  +2.35%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.60'
$ws.Cells.Item(20, 5).Value = '  +1.52%  '
$ws.Cells.Item(21, 5).Value = '  -0.17%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.092'
$ws.Cells.Item(22, 5).Value = '  +3.52%  '
$ws.Cells.Item(23, 4).Value = '28.372.68'
$ws.Cells.Item(23, 5).Value = '  +1.31%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.19'
$ws.Cells.Item(24, 5).Value = '  +1.22%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.172'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.466'
$ws.Cells.Item(26, 5).Value = '  +8.35%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '157.23'
$ws.Cells.Item(27, 5).Value = '  -0.34%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '20.56'
$ws.Cells.Item(28, 5).Value = '  +1.67%  '
$ws.Cells.Item(29, 4).Value = '2.041.42'
$ws.Cells.Item(29, 5).Value = '  +3.14%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '124.11'
$ws.Cells.Item(30, 5).Value = '  +3.74%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.120'
$ws.Cells.Item(31, 5).Value = '  +2.20%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.1093'
$ws.Cells.Item(32, 5).Value = '  +5.01%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.677'
$ws.Cells.Item(33, 5).Value = '  +3.57%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.631'
$ws.Cells.Item(34, 5).Value = '  -0.57%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.07175'
$ws.Cells.Item(35, 5).Value = '  +13.09%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.2255'
$ws.Cells.Item(36, 5).Value = '  +1.11%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.02339'
$ws.Cells.Item(37, 5).Value = '  +3.55%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '5.211'
$ws.Cells.Item(38, 5).Value = '  +5.13%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '8.830'
$ws.Cells.Item(39, 5).Value = '  +5.02%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.6268'
$ws.Cells.Item(40, 5).Value = '  +2.03%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '11.30'
$ws.Cells.Item(41, 5).Value = '  +3.04%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.187'
$ws.Cells.Item(42, 5).Value = '  +1.37%  '
$ws.Cells.Item(43, 5).Value = '  -0.20%  '
$ws.Cells.Item(44, 5).Value = '  -2.47%  '
$ws.Cells.Item(45, 5).Value = '  +1.56%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.706'
$ws.Cells.Item(46, 5).Value = '  +1.17%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.5848'
$ws.Cells.Item(47, 5).Value = '  +1.93%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '125.62'
$ws.Cells.Item(48, 5).Value = '  +0.20%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.986'
$ws.Cells.Item(49, 5).Value = '  +3.46%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.207'
$ws.Cells.Item(50, 5).Value = '  +1.04%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06893'
